$wb = $excel.ActiveWorkbook

# "展览" (Exhibitions) sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 494
$wsExpo.Range("F6").Value = 676

# "全部类型" (All types) sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 494
$wsAll.Range("F6").Value = 676
